# "Purpose List" workbook update — calapan and latest updates
#
# The worksheet holds a single alphabetically-sorted list of "Purpose"
# strings in column B (rows 3..N), numbered sequentially in column A.
# This change:
#   1. Renames one purpose (typo/ID fix: SEL 48713 -> SEL 487B)
#   2. Removes one duplicate purpose entry
#   3. Adds 23 new purpose entries
#   4. Re-sorts the whole list alphabetically (case-insensitive) and
#      renumbers column A sequentially, extending the sheet as needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Read the current purpose list (B3:B338) -------------------------
$firstRow = 3
$lastRowBefore = 338

$srcRange = $ws.Range("B$firstRow`:B$lastRowBefore")
$srcVals = $srcRange.Value2

$purposes = @()
for ($i = $srcVals.GetLowerBound(0); $i -le $srcVals.GetUpperBound(0); $i++) {
    $purposes += $srcVals[$i, 1]
}

# ---- 2. Rename + remove --------------------------------------------------
$renameFrom = "Enclosure and Termination of SEL 48713 Bus Protection Relay"
$renameTo = "Enclosure and Termination of SEL 487B Bus Protection Relay"
$removeText = "Consumables, Tools and Equipment for Spare Stator Rewinding."

$updated = @()
foreach ($p in $purposes) {
    if ($p -eq $renameFrom) {
        $updated += $renameTo
    } elseif ($p -eq $removeText) {
        continue
    } else {
        $updated += $p
    }
}

# ---- 3. New purposes to add ----------------------------------------------
$additions = @(
    "Air Cooler Installation",
    "Construction of Canal in Right Bank of Engine Foundation",
    "Consumables for Cleaning of Foundation Canal",
    "Control Unit & Heating System for Heavy Fuel Oil Purifier",
    "Fabrication of Expansion Belows",
    "Fabrication of Starting/Relief Valve Blind Plugs",
    "Fabrication of Wooden Crates for Engine Parts (Pistons & Connecting Rods) PROGEN",
    "For Jacking Bolt",
    "Governor Maintenance Tools",
    "Inspection of Generator Stator Frame",
    "Installation in Engine Drain Line",
    "Installation of Canopy at 750kVA Transformer Protection Shed",
    "Installation of Centrifugal Pump at Tank Farm Sump Pit",
    "Installation of SEL 487B Relay",
    "Replacement of Air Cooler O-Ring",
    "Replacement of Hydraulic Tensioning Jack Damaged Parts and Additional Hose",
    "Servicing, Reconditioning and Automation of One (1) Unit Westfalia Lube Oil Purfier",
    "Sludge Tank Recovery, Common.",
    "Spare Stator Rewinding Enclosure",
    "Special Tool, Pielstick",
    "Stairway Modification and Construction of Canopy.",
    "Tarpaulin Enclosure Hole Patching",
    "Tools, Inventory-Mary Grace Bugna"
)
$updated += $additions

# ---- 4. Sort alphabetically (case-insensitive, PowerShell default) ------
$sorted = @($updated | Sort-Object)
$n = $sorted.Count
$lastRowAfter = $firstRow + $n - 1

# ---- 5. Extend cell formatting (borders/alignment) for the new rows -----
if ($lastRowAfter -gt $lastRowBefore) {
    $fmtSrc = $ws.Range("A$lastRowBefore`:B$lastRowBefore")
    $fmtDst = $ws.Range("A$($lastRowBefore + 1):B$lastRowAfter")
    [void]$fmtSrc.Copy($fmtDst)
}

# ---- 6. Write the renumbered, re-sorted list back ------------------------
$data = New-Object 'object[,]' $n, 2
for ($i = 0; $i -lt $n; $i++) {
    $data[$i, 0] = $i + 1
    $data[$i, 1] = $sorted[$i]
}
$ws.Range("A$firstRow`:B$lastRowAfter").Value = $data

# ---- 7. Match the final selection shown in the saved workbook -----------
[void]$ws.Range("B$lastRowAfter").Select()
